$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# "simplificamos el modelo del caso 1, habiamos creado una tabla demas"
# The 4th diagram box ("archivos drive", header H4 + rows H5:H7, plus
# the "1  0..n" relationship label in G5) is being removed, and the
# white canvas rectangle that used to just frame that little table is
# redrawn bigger (G2:I10) to cover the area where it used to sit.
# -----------------------------------------------------------------

# Remove the "archivos drive" table entirely (header + its 3 field rows)
$ws.Range("H4:H7").Clear()
# Remove the relationship label that pointed at it
$ws.Range("G5").Clear()

# Re-paint a blank white canvas rectangle (same look as the existing
# G6 filler cell) over the area, extending it out to G2:I10
$ws.Range("G6").Copy()
$ws.Range("G2:I10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the final selection left behind in the saved file
$ws.Range("H18").Select()

Write-Output "done"
